# Applies the two textual edits from the commit:
#   1. Slide 14 ("Rectangle 5"): merge the "Xavier " / "Scherrer :" runs
#      of the first paragraph into a single run "Xavier Scherrer :".
#   2. Slide 15 ("ZoneTexte 5"): split the single run
#      "Amélioration de notre méthode programmation" (3rd paragraph)
#      into two runs: "Amélioration de notre méthode " and
#      "de programmation".

$p = $ppt.ActivePresentation

# -----------------------------------------------------------------
# 1) Slide 14 - merge "Xavier " + "Scherrer :" into one run.
# -----------------------------------------------------------------
$slide14 = $p.Slides.Item(14)
$shape14 = $slide14.Shapes.Item(3)            # "Rectangle 5"
$tr14 = $shape14.TextFrame.TextRange
$para1 = $tr14.Paragraphs(1, 1)

$run1 = $para1.Runs(1, 1)
$run2 = $para1.Runs(2, 1)

# Put the full merged text on the first run, then clear out the
# second run's text so the (now redundant) run node is dropped.
$run1.Text = "Xavier Scherrer :"
$para1Fresh = $tr14.Paragraphs(1, 1)
$run2Fresh = $para1Fresh.Runs(2, 1)
$run2Fresh.Text = ""

# -----------------------------------------------------------------
# 2) Slide 15 - split "Amélioration de notre méthode programmation"
#    into "Amélioration de notre méthode " + "de programmation".
# -----------------------------------------------------------------
$slide15 = $p.Slides.Item(15)
$shape15 = $slide15.Shapes.Item(3)            # "ZoneTexte 5"
$tr15 = $shape15.TextFrame.TextRange
$para3 = $tr15.Paragraphs(3, 1)

$runA = $para3.Runs(1, 1)
$runA.Text = "Amélioration de notre méthode "
$runB = $runA.InsertAfter("de programmation")

Write-Output "Edits applied."
